# Append the latest day's profit row (run date 2025-08-21) to Sheet1.
#
# A4 must stay a literal text value "08/21/2025" (matching the existing
# date-as-text cells in A2/A3), not get auto-converted into a serial
# date number. Forcing NumberFormat to "@" (Text) before assigning the
# value prevents Excel's date auto-detection; resetting the style back
# to "Normal" afterwards avoids leaving a stray text-format style on
# the cell so the sheet's styling stays identical to the other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "08/21/2025"
$ws.Cells.Item(4, 1).Style = "Normal"

$ws.Cells.Item(4, 2).Value = 13491.24
